$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-4.01%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'30.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.94%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.940"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.37%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07170"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-8.41%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.788"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-11.57%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.656"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.23%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.744"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.97%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.8945"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.77%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1649"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-6.33%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07706"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.26%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08052"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-6.92%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03059"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-3.32%"
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'-0.11%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001507"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.52%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005813"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.76%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.471"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.23%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E19").Value = "'-0.92%"
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'-1.32%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.042"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-5.57%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.1999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.43%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04515"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.18%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001213"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.97%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'-9.99%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'-0.13%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.01599"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-8.20%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04380"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-8.43%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007336"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.69%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1308"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-4.19%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.007656"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.002049"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-13.25%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.009254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-12.55%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00005951"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.99%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.246"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'172.73%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.002999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-3.30%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'-0.08%"
$ws.Range("E51").Style = "Normal"
